# regen sval data to filter save games
# Updates the numeric stat columns (B:G) for rows 2-8 on the active sheet
# to reflect the re-computed values after filtering save games.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.505614041169197
$ws.Range("C2").Value = 0.3375848360084654
$ws.Range("D2").Value = 3.082599426703578
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 11.40722638115186

$ws.Range("B3").Value = 1.505614041169197
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 3.082599426703578
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 6.741336633845642

$ws.Range("B4").Value = 3.182878228561681
$ws.Range("C4").Value = 1.65323645889881
$ws.Range("D4").Value = 0.1529057820181812
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 5.488907176552729

$ws.Range("B5").Value = 1.505614041169197
$ws.Range("C5").Value = 1.65323645889881
$ws.Range("D5").Value = 0.1529057820181812
$ws.Range("E5").Value = 0.4998867070740569
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 3.811642989160245

$ws.Range("B6").Value = 3.182878228561681
$ws.Range("C6").Value = 1.65323645889881
$ws.Range("D6").Value = 3.082599426703578
$ws.Range("E6").Value = 246.9852506941017
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 254.9039648082657

$ws.Range("B7").Value = 0.1554434735375247
$ws.Range("C7").Value = 0.3375848360084654
$ws.Range("D7").Value = 157.8057217802531
$ws.Range("E7").Value = 6.48142807727062
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 164.7801781670697

$ws.Range("B8").Value = 0.3464964993005633
$ws.Range("C8").Value = 1.65323645889881
$ws.Range("D8").Value = 16.98373111632243
$ws.Range("E8").Value = 6.48142807727062
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 25.46489215179242
